$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells that look numeric stay as text, matching the source data format
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '43.475.56'
$ws.Cells.Item(2, 5).Value = '  +3.21%  '

$ws.Cells.Item(3, 4).Value = '2.331.99'
$ws.Cells.Item(3, 5).Value = '  +3.37%  '

$ws.Cells.Item(4, 5).Value = '  +0.10%  '

$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Cells.Item(5, 4).Value = '309.13'
$ws.Cells.Item(5, 5).Value = '  +0.73%  '

$ws.Cells.Item(6, 2).Value = 'Solana'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(6, 4).Value = '105.51'
$ws.Cells.Item(6, 5).Value = '  +7.71%  '

$ws.Cells.Item(7, 5).Value = '  +0.62%  '

$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 4).Value = '0.522'
$ws.Cells.Item(9, 5).Value = '  +6.49%  '

$ws.Cells.Item(10, 4).Value = '36.36'
$ws.Cells.Item(10, 5).Value = '  +2.76%  '

$ws.Cells.Item(11, 4).Value = '52.98'

$ws.Cells.Item(12, 5).Value = '  -0.13%  '

$ws.Cells.Item(13, 5).Value = '  -1.55%  '

$ws.Cells.Item(14, 4).Value = '7.01'
$ws.Cells.Item(14, 5).Value = '  +3.83%  '

$ws.Cells.Item(15, 4).Value = '2.683.95'
$ws.Cells.Item(15, 5).Value = '  +3.23%  '

$ws.Cells.Item(16, 4).Value = '15.46'
$ws.Cells.Item(16, 5).Value = '  +6.60%  '

$ws.Cells.Item(17, 4).Value = '2.320.21'
$ws.Cells.Item(17, 5).Value = '  +2.89%  '

$ws.Cells.Item(18, 4).Value = '0.805'
$ws.Cells.Item(18, 5).Value = '  +3.21%  '

$ws.Cells.Item(19, 4).Value = '43.444.39'
$ws.Cells.Item(19, 5).Value = '  +3.42%  '

$ws.Cells.Item(20, 4).Value = '11.97'
$ws.Cells.Item(20, 5).Value = '  -1.88%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0925'
$ws.Cells.Item(21, 5).Value = '  +2.20%  '

$ws.Cells.Item(22, 4).Value = '6.25'
$ws.Cells.Item(22, 5).Value = '  +5.11%  '

$ws.Cells.Item(23, 4).Value = '68.20'
$ws.Cells.Item(23, 5).Value = '  +1.62%  '

$ws.Cells.Item(24, 4).Value = '241.96'
$ws.Cells.Item(24, 5).Value = '  +2.41%  '

$ws.Cells.Item(25, 4).Value = '2.04'
$ws.Cells.Item(25, 5).Value = '  +4.77%  '

$ws.Cells.Item(26, 5).Value = '  +1.49%  '

$ws.Cells.Item(27, 5).Value = '  +0.15%  '

$ws.Cells.Item(28, 4).Value = '25.04'
$ws.Cells.Item(28, 5).Value = '  +7.13%  '

$ws.Cells.Item(29, 4).Value = '36.50'
$ws.Cells.Item(29, 5).Value = '  -4.22%  '

$ws.Cells.Item(30, 5).Value = '  +3.66%  '

$ws.Cells.Item(31, 5).Value = '  +1.20%  '

$ws.Cells.Item(32, 4).Value = '163.03'
$ws.Cells.Item(32, 5).Value = '  -2.10%  '

$ws.Cells.Item(33, 4).Value = '5.26'
$ws.Cells.Item(33, 5).Value = '  +1.70%  '

$ws.Cells.Item(34, 5).Value = '  +0.11%  '

$ws.Cells.Item(35, 4).Value = '18.23'
$ws.Cells.Item(35, 5).Value = '  +3.62%  '

$ws.Cells.Item(36, 5).Value = '  +6.76%  '

$ws.Cells.Item(37, 5).Value = '  +2.44%  '

$ws.Cells.Item(38, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(38, 4).Value = '3.06'
$ws.Cells.Item(38, 5).Value = '  +0.47%  '

$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(39, 4).Value = '4.61'
$ws.Cells.Item(39, 5).Value = '  +12.48%  '

$ws.Cells.Item(40, 2).Value = 'ARBITRUM'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(40, 4).Value = '1.88'
$ws.Cells.Item(40, 5).Value = '  +5.27%  '

$ws.Cells.Item(41, 2).Value = 'Kaspa'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(41, 4).Value = '0.107'
$ws.Cells.Item(41, 5).Value = '  +3.98%  '

$ws.Cells.Item(42, 5).Value = '  +0.55%  '

$ws.Cells.Item(43, 4).Value = '2.47'
$ws.Cells.Item(43, 5).Value = '  +11.14%  '

$ws.Cells.Item(44, 5).Value = '  +3.29%  '

$ws.Cells.Item(45, 4).Value = '1.973.85'
$ws.Cells.Item(45, 5).Value = '  +1.95%  '

$ws.Cells.Item(46, 4).Value = '19.21'
$ws.Cells.Item(46, 5).Value = '  +3.10%  '

$ws.Cells.Item(47, 4).Value = '3.07'
$ws.Cells.Item(47, 5).Value = '  +6.30%  '

$ws.Cells.Item(48, 5).Value = '  +6.49%  '

$ws.Cells.Item(49, 4).Value = '58.22'
$ws.Cells.Item(49, 5).Value = '  +8.05%  '

$ws.Cells.Item(50, 4).Value = '2.93'
$ws.Cells.Item(50, 5).Value = '  +0.30%  '

$ws.Cells.Item(51, 4).Value = '1.59'
$ws.Cells.Item(51, 5).Value = '  +8.24%  '
